$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

# New Gasoline_Price (col B), Diesel_Price (col C), LPG_Price (col D) values
# per country row (rows 2-27), refreshed from the underlying Power Query data source.

$gas = @(
    "184,78 ", "139,20 ", "141,41 ", "210,52 ", "165,63 ", "205,07 ",
    "182,32 ", "177,51 ", "205,53 ", "222,53 ", "203,19 ", "166,49 ",
    "170,01 ", "199,32 ", "209,64 ", "161,57 ", "179,04 ", "141,94 ",
    "236,84 ", "160,31 ", "206,47 ", "163,90 ", "176,62 ", "171,30 ",
    "179,97 ", "125,70 "
)

$diesel = @(
    "182,55 ", "135,60 ", "141,41 ", "226,28 ", "155,42 ", "184,66 ",
    "164,26 ", "166,96 ", "189,59 ", "226,29 ", "170,01 ", "151,25 ",
    "163,56 ", "192,40 ", "197,68 ", "163,68 ", "170,24 ", "125,74 ",
    "203,42 ", "157,56 ", "190,64 ", "165,30 ", "186,72 ", "171,65 ",
    "167,55 ", "124,26 "
)

$lpg = @(
    "146,32 ", "79,80 ", "71,30 ", "116,70 ", "79,20 ", "128,03 ",
    "105,29 ", "111,74 ", "119,01 ", "126,23 ", "102,71 ", "85,59 ",
    "97,99 ", "117,25 ", "82,42 ", "102,94 ", "105,99 ", "80,02 ",
    "93,09 ", "81,25 ", "104,23 ", "82,76 ", "94,95 ", "97,78 ",
    "83,36 ", "66,58 "
)

for ($i = 0; $i -lt 26; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $gas[$i]
    $ws.Range("C$row").Value = $diesel[$i]
    $ws.Range("D$row").Value = $lpg[$i]
}
